# Fix bug trong file convertJsonToExcel
# Insert 4 new customer rows right after the header row (row 1), pushing
# all existing data down by 4 rows, and populate the new rows with the
# new customers' data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at rows 2-5 (existing data shifts down to rows 6-139)
$ws.Range("A2:J5").EntireRow.Insert()

# Row 2 - sơn thị chành thi
$ws.Cells.Item(2, 1).Value = "KH"
$ws.Cells.Item(2, 2).Value = 438
$ws.Cells.Item(2, 3).Value = "sơn thị chành thi"
$ws.Cells.Item(2, 4).Value = "SÓC TRĂNG"
$ws.Cells.Item(2, 6).Value = "'0327114485"
$ws.Cells.Item(2, 9).Value = 4000000
$ws.Cells.Item(2, 10).Value = 0

# Row 3 - triệu tú kiều
$ws.Cells.Item(3, 1).Value = "KH"
$ws.Cells.Item(3, 2).Value = 437
$ws.Cells.Item(3, 3).Value = "triệu tú kiều "
$ws.Cells.Item(3, 4).Value = "SÓC TRĂNG"
$ws.Cells.Item(3, 6).Value = "'0974395268"
$ws.Cells.Item(3, 9).Value = 12000000
$ws.Cells.Item(3, 10).Value = 0

# Row 4 - lý thị thuý vi
$ws.Cells.Item(4, 1).Value = "KH"
$ws.Cells.Item(4, 2).Value = 436
$ws.Cells.Item(4, 3).Value = "lý thị thuý vi"
$ws.Cells.Item(4, 4).Value = "SÓC TRĂNG"
$ws.Cells.Item(4, 6).Value = "'0396202865"
$ws.Cells.Item(4, 9).Value = 5000000
$ws.Cells.Item(4, 10).Value = 0

# Row 5 - nguyễn thị kim phượng
$ws.Cells.Item(5, 1).Value = "KH"
$ws.Cells.Item(5, 2).Value = 435
$ws.Cells.Item(5, 3).Value = "nguyễn thị kim phượng "
$ws.Cells.Item(5, 4).Value = "SÓC TRĂNG"
$ws.Cells.Item(5, 6).Value = "'0368976358"
$ws.Cells.Item(5, 9).Value = 13000000
$ws.Cells.Item(5, 10).Value = 0

Write-Output "done"
